# bot5 finalizado version 1.01
# "parametrosInicio" is the active sheet (activeTab="1" in the source file).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10 ("Validacion") changes from "ASIG, FECHA, IMP, NOMBRE" to "FECHA, IMP, NOMBRE"
$ws.Range("B10").Value = "FECHA, IMP, NOMBRE"

# B9 ("Tipo de cuenta") keeps its displayed text "CUENTA ETV"
$ws.Range("B9").Value = "CUENTA ETV"

# Move/update the sheet selection from B10 to D10
$ws.Range("D10").Select()
